# adicionado politica de preco
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C (so old C,D,E -> E,F,G)
$ws.Range("C1:D1").EntireColumn.Insert()

# New header cells
$ws.Range("C1").Value = "modelo"
$ws.Range("D1").Value = "politica"

# lowercase the "tipo" column (now column F) values
$ws.Range("F2").Value = "classico"
$ws.Range("F3").Value = "classico"
$ws.Range("F4").Value = "classico"
$ws.Range("F5").Value = "classico"
$ws.Range("F6").Value = "classico"
$ws.Range("F7").Value = "premium"
$ws.Range("F8").Value = "classico"
$ws.Range("F9").Value = "classico"
$ws.Range("F10").Value = "premium"
$ws.Range("F11").Value = "classico"
$ws.Range("F12").Value = "premium"
$ws.Range("F13").Value = "premium"
$ws.Range("F14").Value = "premium"
$ws.Range("F15").Value = "premium"

# new "modelo" values
$ws.Range("C2").Value = "FONTE 200 MONO"
$ws.Range("C3").Value = "FONTE 200 BOB"
$ws.Range("C4").Value = "FONTE 200A"
$ws.Range("C5").Value = "FONTE 120 BOB"
$ws.Range("C6").Value = "FONTE 200A LITE"
$ws.Range("C7").Value = "Sem Modelo"
$ws.Range("C8").Value = "FONTE 90 BOB"
$ws.Range("C9").Value = "FONTE 120A"
$ws.Range("C10").Value = "FONTE 120A"
$ws.Range("C11").Value = "FONTE 120A"
$ws.Range("C12").Value = "FONTE 70A LITE"
$ws.Range("C13").Value = "FONTE 200 BOB"
$ws.Range("C14").Value = "FONTE 200A"
$ws.Range("C15").Value = "FONTE 70A LITE"

# new "politica" values
$ws.Range("D2").Value = "Igual"
$ws.Range("D3").Value = "Igual"
$ws.Range("D4").Value = "Igual"
$ws.Range("D5").Value = "Igual"
$ws.Range("D6").Value = "Igual"
$ws.Range("D7").Value = ""
$ws.Range("D8").Value = "Igual"
$ws.Range("D9").Value = "Igual"
$ws.Range("D10").Value = "Igual"
$ws.Range("D11").Value = "Acima"
$ws.Range("D12").Value = "Acima"
$ws.Range("D13").Value = "Igual"
$ws.Range("D14").Value = "Igual"
$ws.Range("D15").Value = "Acima"

# row 14's "full" column (now E14) changes from FULL to NA
$ws.Range("E14").Value = "NA"

# update links (now column G) with the new tracking_id and position values
$ws.Range("G2").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-200a-storm-voltimetro-digital-mono-220v-cor-preto/p/MLB24006449?pdp_filters=seller_id:703818843#searchVariation=MLB24006449&position=5&search_layout=stack&type=product&tracking_id=986a7aa9-7d24-47d8-a753-98e6da131914"
$ws.Range("G3").Value = "https://www.mercadolivre.com.br/fonte-automotiva-jfa-storm-200a-bob-carregador-automatico-bivolt-cor-bob-200a-jfa/p/MLB24834408?pdp_filters=seller_id:703818843#searchVariation=MLB24834408&position=2&search_layout=stack&type=product&tracking_id=986a7aa9-7d24-47d8-a753-98e6da131914"
$ws.Range("G4").Value = "https://www.mercadolivre.com.br/fonte-carregador-automotiva-jfa-200a-slim-bivolt-voltimetro/p/MLB21348561?pdp_filters=seller_id:703818843#searchVariation=MLB21348561&position=8&search_layout=stack&type=product&tracking_id=986a7aa9-7d24-47d8-a753-98e6da131914"
$ws.Range("G5").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-120a-bob-slim-bivolt-cor-preto/p/MLB22144397?pdp_filters=seller_id:703818843#searchVariation=MLB22144397&position=4&search_layout=stack&type=product&tracking_id=986a7aa9-7d24-47d8-a753-98e6da131914"
$ws.Range("G6").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-200a-lite-storm-slim-bivolt-cor-azul/p/MLB24154371?pdp_filters=seller_id:703818843#searchVariation=MLB24154371&position=3&search_layout=stack&type=product&tracking_id=986a7aa9-7d24-47d8-a753-98e6da131914"
$ws.Range("G7").Value = "https://produto.mercadolivre.com.br/MLB-2698193784-controle-longa-distncia-jfa-acqua-1200-resistente-a-agua-_JM#position%3D9%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D986a7aa9-7d24-47d8-a753-98e6da131914"
$ws.Range("G8").Value = "https://www.mercadolivre.com.br/fonte-carregador-jfa-bob-storm-90a-bivolt-automatico-cor-preto/p/MLB21562641?pdp_filters=seller_id:703818843#searchVariation=MLB21562641&position=6&search_layout=stack&type=product&tracking_id=986a7aa9-7d24-47d8-a753-98e6da131914"
$ws.Range("G9").Value = "https://www.mercadolivre.com.br/fonte-automotiva-120a-amperes-jfa-carregador-cor-preto/p/MLB21392652?pdp_filters=seller_id:703818843#searchVariation=MLB21392652&position=1&search_layout=stack&type=product&tracking_id=986a7aa9-7d24-47d8-a753-98e6da131914"
$ws.Range("G10").Value = "https://produto.mercadolivre.com.br/MLB-2164283206-fonte-carregador-jfa-120a-storm-bivolt-com-medidor-cca-_JM#position%3D10%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D986a7aa9-7d24-47d8-a753-98e6da131914"
$ws.Range("G11").Value = "https://produto.mercadolivre.com.br/MLB-2164246906-fonte-carregador-jfa-120a-storm-bivolt-com-medidor-cca-_JM#position%3D11%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D986a7aa9-7d24-47d8-a753-98e6da131914"
$ws.Range("G12").Value = "https://produto.mercadolivre.com.br/MLB-3801613892-fonte-carregador-de-bateria-jfa-70a-lite-storm-slim-bivolt-_JM#position%3D12%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D986a7aa9-7d24-47d8-a753-98e6da131914"
$ws.Range("G13").Value = "https://produto.mercadolivre.com.br/MLB-2774547578-fonte-carregadora-jfa-bob-storm-200a-bivolt-_JM#position%3D13%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D986a7aa9-7d24-47d8-a753-98e6da131914"
$ws.Range("G14").Value = "https://produto.mercadolivre.com.br/MLB-2677379815-fonte-automotiva-jfa-200-amperes-storm-bivolt-c-medidor-caa-_JM#position%3D14%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D986a7aa9-7d24-47d8-a753-98e6da131914"
$ws.Range("G15").Value = "https://produto.mercadolivre.com.br/MLB-3813696866-fonte-carregador-de-bateria-jfa-70a-lite-storm-slim-bivolt-_JM#position%3D15%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3D986a7aa9-7d24-47d8-a753-98e6da131914"
